$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.607.30"
$ws.Range("E2").Value = "  -3.00%  "
$ws.Range("D3").Value = "1.981.89"
$ws.Range("E3").Value = "  -3.79%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'246.48"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").Value = "'0.635"
$ws.Range("E6").Value = "  -4.57%  "
$ws.Range("D7").Value = "'58.09"
$ws.Range("E7").Value = "  +5.84%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'58.69"
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("D10").Value = "'0.361"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").Value = "'0.0736"
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("D13").Value = "'0.959"
$ws.Range("E13").Value = "  +2.68%  "
$ws.Range("D14").Value = "'14.58"
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("D15").Value = "2.271.80"
$ws.Range("E15").Value = "  -3.70%  "
$ws.Range("D16").Value = "'5.29"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").Value = "1.979.20"
$ws.Range("E17").Value = "  -3.69%  "
$ws.Range("D18").Value = "'18.46"
$ws.Range("E18").Value = "  +7.30%  "
$ws.Range("D19").Value = "35.571.05"
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("D20").Value = "'71.45"
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("D21").Value = "0.0₃0848"
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").Value = "'232.98"
$ws.Range("E23").Value = "  -2.25%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'2.60"
$ws.Range("E25").Value = "  +21.36%  "
$ws.Range("D26").Value = "'2.28"
$ws.Range("E26").Value = "  -3.98%  "
$ws.Range("D27").Value = "'164.97"
$ws.Range("E28").Value = "  -1.95%  "
$ws.Range("D29").Value = "'19.23"
$ws.Range("E29").Value = "  -4.73%  "
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("D31").Value = "'4.88"
$ws.Range("E31").Value = "  -4.11%  "
$ws.Range("E32").Value = "  -7.88%  "
$ws.Range("D33").Value = "'0.0954"
$ws.Range("E33").Value = "  +13.16%  "
$ws.Range("D34").Value = "'0.0595"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("E35").Value = "  +9.31%  "
$ws.Range("D36").Value = "'4.37"
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "'1.77"
$ws.Range("E38").Value = "  -3.41%  "
$ws.Range("D39").Value = "'5.51"
$ws.Range("E39").Value = "  +10.90%  "
$ws.Range("D40").Value = "'1.24"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.09"
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'7.79"
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("D45").Value = "'93.63"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").Value = "'0.0904"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("D48").Value = "1.374.49"
$ws.Range("E48").Value = "  -2.74%  "
$ws.Range("D49").Value = "'2.89"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").Value = "'46.93"
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("D51").Value = "'2.27"
$ws.Range("E51").Value = "  -0.04%  "
